$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1666.9375
$ws.Range("I80").Value = 2424.5715
$ws.Range("J80").Value = 1077.6666
$ws.Range("K80").Value = 7273.7145
$ws.Range("L80").Value = 3232.9998
$ws.Range("M80").Value = -6275.7145
$ws.Range("N80").Value = -5228.9998
$ws.Range("H83").Value = 1666.9375
$ws.Range("I83").Value = 2424.5715
$ws.Range("J83").Value = 1077.6666
$ws.Range("K83").Value = 21821.1435
$ws.Range("L83").Value = 9698.999400000001
$ws.Range("M83").Value = -16829.1435
$ws.Range("N83").Value = -19682.9994
$ws.Range("H86").Value = 1487.125
$ws.Range("J86").Value = 1382.8334
$ws.Range("L86").Value = 1382.8334
$ws.Range("N86").Value = -3628.8334
$ws.Range("H89").Value = 1487.125
$ws.Range("J89").Value = 1382.8334
$ws.Range("L89").Value = 6914.166999999999
$ws.Range("N89").Value = -18146.167
$ws.Range("H98").Value = 1414.7084
$ws.Range("I98").Value = 1193.3
$ws.Range("J98").Value = 2521.75
$ws.Range("K98").Value = 1193.3
$ws.Range("L98").Value = 2521.75
$ws.Range("M98").Value = 304.7
$ws.Range("N98").Value = -5517.75
$ws.Range("H122").Value = 1414.7084
$ws.Range("I122").Value = 1193.3
$ws.Range("J122").Value = 2521.75
$ws.Range("K122").Value = 3579.9
$ws.Range("L122").Value = 7565.25
$ws.Range("M122").Value = -1129.9
$ws.Range("N122").Value = -12465.25
$ws.Range("H123").Value = 43333.332
$ws.Range("J123").Value = 43333.332
$ws.Range("L123").Value = 43333.332
$ws.Range("N123").Value = -53133.332
$ws.Range("H132").Value = 999.3125
$ws.Range("I132").Value = 912.1667
$ws.Range("J132").Value = 1260.75
$ws.Range("K132").Value = 2736.5001
$ws.Range("L132").Value = 3782.25
$ws.Range("M132").Value = -206.5001000000002
$ws.Range("N132").Value = -8842.25
$ws.Range("H138").Value = 2517.7678
$ws.Range("J138").Value = 2479.64
$ws.Range("L138").Value = 7438.92
$ws.Range("N138").Value = -17718.92
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4599.268
$ws.Range("I32").Value = 3274.6538
$ws.Range("J32").Value = 21819.25
$ws.Range("K32").Value = 3274.6538
$ws.Range("L32").Value = 21819.25
$ws.Range("M32").Value = -2987.6538
$ws.Range("N32").Value = -22393.25
$ws.Range("H74").Value = 1277.5927
$ws.Range("I74").Value = 489.95
$ws.Range("J74").Value = 3528
$ws.Range("K74").Value = 489.95
$ws.Range("L74").Value = 3528
$ws.Range("M74").Value = 384.05
$ws.Range("N74").Value = -5276
$ws.Range("H77").Value = 1277.5927
$ws.Range("I77").Value = 489.95
$ws.Range("J77").Value = 3528
$ws.Range("K77").Value = 2449.75
$ws.Range("L77").Value = 17640
$ws.Range("M77").Value = 1918.25
$ws.Range("N77").Value = -26376
$ws.Range("H132").Value = 1960.3334
$ws.Range("I132").Value = 1531.6471
$ws.Range("J132").Value = 3001.4285
$ws.Range("K132").Value = 4594.9413
$ws.Range("L132").Value = 9004.2855
$ws.Range("M132").Value = -2064.9413
$ws.Range("N132").Value = -14064.2855
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M68").Value = -46189
$ws.Range("H68").Value = 32125
$ws.Range("I68").Value = 47000
$ws.Range("J68").Value = 30000
$ws.Range("K68").Value = 47000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622
$ws.Range("M71").Value = -136944
$ws.Range("H71").Value = 32125
$ws.Range("I71").Value = 47000
$ws.Range("J71").Value = 30000
$ws.Range("K71").Value = 141000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112
$ws.Range("I86").Value = 9600
$ws.Range("J86").Value = 502252.5
$ws.Range("K86").Value = 9600
$ws.Range("L86").Value = 502252.5
$ws.Range("M86").Value = -8477
$ws.Range("N86").Value = -504498.5
$ws.Range("I89").Value = 9600
$ws.Range("J89").Value = 502252.5
$ws.Range("K89").Value = 48000
$ws.Range("L89").Value = 2511262.5
$ws.Range("M89").Value = -42384
$ws.Range("N89").Value = -2522494.5
$ws.Range("H134").Value = 6007.16
$ws.Range("I134").Value = 6752.8
$ws.Range("K134").Value = 20258.4
$ws.Range("M134").Value = -17723.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 400.2
$ws.Range("I7").Value = 467
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 467
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -354
$ws.Range("N7").Value = -526
$ws.Range("H62").Value = 3974.75
$ws.Range("J62").Value = 4450
$ws.Range("L62").Value = 4450
$ws.Range("N62").Value = -5698
$ws.Range("H65").Value = 3974.75
$ws.Range("J65").Value = 4450
$ws.Range("L65").Value = 22250
$ws.Range("N65").Value = -28490
$ws.Range("H92").Value = 41249.25
$ws.Range("J92").Value = 41249.25
$ws.Range("L92").Value = 41249.25
$ws.Range("N92").Value = -46241.25
$ws.Range("H132").Value = 2342.2222
$ws.Range("I132").Value = 1561.6
$ws.Range("K132").Value = 4684.799999999999
$ws.Range("M132").Value = -2154.799999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N92").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H97").Value = 750.25
$ws.Range("I97").Value = 399.66666
$ws.Range("K97").Value = 1198.99998
$ws.Range("M97").Value = -702.9999800000001
$ws.Range("N136").Value = -15600
$ws.Range("H136").Value = 1129.3334
$ws.Range("I136").Value = 995.2
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 2985.6
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = 2114.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N80").Value = -5996
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 4000
$ws.Range("N83").Value = -29984
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 20000
$ws.Range("H102").Value = 4023.4
$ws.Range("I102").Value = 4874.5
$ws.Range("K102").Value = 4874.5
$ws.Range("M102").Value = -3252.5
$ws.Range("H122").Value = 1682.5
$ws.Range("I122").Value = 1618.8667
$ws.Range("K122").Value = 4856.6001
$ws.Range("M122").Value = -2406.6001
$ws.Range("H126").Value = 2573269.8
$ws.Range("I126").Value = 4632320.5
$ws.Range("K126").Value = 13896961.5
$ws.Range("M126").Value = -13894491.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1865.6562
$ws.Range("I7").Value = 1809.0344
$ws.Range("K7").Value = 1809.0344
$ws.Range("M7").Value = -1697.0344
$ws.Range("H40").Value = 9947.736999999999
$ws.Range("I40").Value = 10478
$ws.Range("J40").Value = 8798.833000000001
$ws.Range("K40").Value = 10478
$ws.Range("L40").Value = 8798.833000000001
$ws.Range("M40").Value = -10342
$ws.Range("N40").Value = -9070.833000000001
$ws.Range("N100").Value = -3642
$ws.Range("H100").Value = 1745.3334
$ws.Range("I100").Value = 930.6667
$ws.Range("J100").Value = 2560
$ws.Range("K100").Value = 930.6667
$ws.Range("L100").Value = 2560
$ws.Range("M100").Value = -389.6667
$ws.Range("H122").Value = 11562.625
$ws.Range("I122").Value = 11100.2
$ws.Range("K122").Value = 33300.60000000001
$ws.Range("M122").Value = -30850.60000000001
$ws.Range("H126").Value = 1865.6562
$ws.Range("I126").Value = 1809.0344
$ws.Range("K126").Value = 5427.1032
$ws.Range("M126").Value = -2957.1032
$ws.Range("H132").Value = 1563
$ws.Range("J132").Value = 1853.3914
$ws.Range("L132").Value = 5560.174199999999
$ws.Range("N132").Value = -10620.1742
$ws.Range("H136").Value = 3242
$ws.Range("I136").Value = 3358
$ws.Range("K136").Value = 10074
$ws.Range("M136").Value = -7524
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N103").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("H122").Value = 30842.77
$ws.Range("I122").Value = 34591.914
$ws.Range("J122").Value = 2099.3333
$ws.Range("K122").Value = 103775.742
$ws.Range("L122").Value = 6297.999899999999
$ws.Range("M122").Value = -101325.742
$ws.Range("N122").Value = -11197.9999
$ws.Range("H136").Value = 1813.579
$ws.Range("I136").Value = 1480.1
$ws.Range("J136").Value = 2184.111
$ws.Range("K136").Value = 4440.299999999999
$ws.Range("L136").Value = 6552.333
$ws.Range("M136").Value = -1890.299999999999
$ws.Range("N136").Value = -11652.333
